# Apply weekly update: insert a new data row at row 44 (pushing existing
# rows 44-93 down to 45-94) and populate the new row with the latest
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44; Excel shifts rows 44:93 down to 45:94
# and copies formatting (including the date-style on column D) from the
# row above, same as a normal Excel "Insert Sheet Rows" operation.
$ws.Rows("44:44").Insert()

# Populate the newly inserted row 44 with the new weekly observation.
$ws.Cells.Item(44, 1).Value = 7
$ws.Cells.Item(44, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(44, 3).Value = "Ñuble"
$ws.Cells.Item(44, 4).Value = 45225
$ws.Cells.Item(44, 5).Value = 16
$ws.Cells.Item(44, 6).Value = 100112026
$ws.Cells.Item(44, 7).Value = "Haba"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 30
$ws.Cells.Item(44, 11).Value = 12000
$ws.Cells.Item(44, 12).Value = 12000
$ws.Cells.Item(44, 13).Value = 12000
$ws.Cells.Item(44, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(44, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(44, 16).Value = 480
$ws.Cells.Item(44, 17).Value = 25
$ws.Cells.Item(44, 18).Value = "Hortaliza"
